# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (interest-count) figures in column F, across
# the 展览 / 演出 / 本地生活 / 全部类型 sheets, to match the latest scrape.

$wb = $excel.ActiveWorkbook

# 展览
$ws = $wb.Worksheets("展览")
$ws.Range("F4").Value  = 7868
$ws.Range("F5").Value  = 7868
$ws.Range("F9").Value  = 8640
$ws.Range("F10").Value = 8640
$ws.Range("F14").Value = 5775
$ws.Range("F16").Value = 2764
$ws.Range("F21").Value = 618
$ws.Range("F22").Value = 77
$ws.Range("F23").Value = 3923
$ws.Range("F28").Value = 181
$ws.Range("F29").Value = 26
$ws.Range("F30").Value = 5481
$ws.Range("F31").Value = 8
$ws.Range("F32").Value = 72
$ws.Range("F34").Value = 392
$ws.Range("F35").Value = 161
$ws.Range("F37").Value = 2529
$ws.Range("F38").Value = 1520
$ws.Range("F40").Value = 1122
$ws.Range("F41").Value = 4871
$ws.Range("F42").Value = 78
$ws.Range("F45").Value = 3583
$ws.Range("F51").Value = 20

# 演出
$ws = $wb.Worksheets("演出")
$ws.Range("F3").Value = 152

# 本地生活
$ws = $wb.Worksheets("本地生活")
$ws.Range("F3").Value = 1367

# 全部类型
$ws = $wb.Worksheets("全部类型")
$ws.Range("F3").Value  = 1367
$ws.Range("F5").Value  = 7868
$ws.Range("F6").Value  = 7868
$ws.Range("F9").Value  = 8640
$ws.Range("F10").Value = 8640
$ws.Range("F13").Value = 5775
$ws.Range("F15").Value = 2764
$ws.Range("F20").Value = 152
$ws.Range("F21").Value = 618
$ws.Range("F22").Value = 77
$ws.Range("F23").Value = 3923
$ws.Range("F28").Value = 181
$ws.Range("F29").Value = 26
$ws.Range("F30").Value = 5481
$ws.Range("F31").Value = 8
$ws.Range("F32").Value = 72
$ws.Range("F33").Value = 392
$ws.Range("F34").Value = 161
$ws.Range("F37").Value = 2529
$ws.Range("F38").Value = 1520
$ws.Range("F40").Value = 1122
$ws.Range("F42").Value = 4871
$ws.Range("F43").Value = 78
$ws.Range("F46").Value = 3583
$ws.Range("F49").Value = 20
